# GestiBank/PATHS.xlsx - "merge eclipse (on y croit)"
# - Remove the old "GET <site>/notifications/client/{id}" row (row 20)
# - Replace it with three new notification-related endpoints appended
#   after the existing table rows (rows 29-31, leaving row 28 blank like
#   the pre-existing gaps at rows 3 and 28)
# - Grow the Tableau1 table/autofilter range to cover the new rows
# - Leave the selection on B37 (next free cell below the table)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Remove the old "getNotifications(clientID)" row entirely.
$ws.Range("A20:E20").ClearContents() | Out-Null

# 2. Append the new notification endpoints.
$ws.Range("A29").Value = "DELETE"
$ws.Range("B29").Value = "<site>/notifications/{id}"
$ws.Range("C29").Value = "deleteNotification(id)"
$ws.Range("D29").Value = "void"
$ws.Range("E29").Value = "CLIENT"

$ws.Range("A30").Value = "POST"
$ws.Range("B30").Value = "<site>/notifications/{id}"
$ws.Range("C30").Value = "addNotificationToClient(clientId)"
$ws.Range("D30").Value = "void"
$ws.Range("E30").Value = "CLIENT"

$ws.Range("A31").Value = "GET"
$ws.Range("B31").Value = "<site>/notifications/{id}"
$ws.Range("C31").Value = "getNotificationsByClient(clientID)"
$ws.Range("D31").Value = "List<Notification>"
$ws.Range("E31").Value = "CLIENT"

# 3. Grow the table (and its autofilter) to the new extent.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E31")) | Out-Null

# 4. Park the selection where the author left it.
$ws.Range("B37").Select() | Out-Null
